# Atualização de bases das ligas, do dia: 16-05-2024 às 23:38
# Swap the data of row 83 (id 81) and row 84 (id 82): these two match
# records ("Forge FC vs Atletico Ottawa" / "Cavalry FC vs Pacific FC CA")
# had been mapped to the wrong fixture id/odds; this restores the correct
# pairing by exchanging every data column (B, E:AB) between the two rows
# while leaving A (row index), C (Div) and D (Date) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 83 (previously row 84's data)
$ws.Range("B83").Value = 6227884
$ws.Range("E83").Value = "Cavalry FC"
$ws.Range("F83").Value = "Pacific FC CA"
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = "H"
$ws.Range("J83").Value = 2.25
$ws.Range("K83").Value = 3.1
$ws.Range("L83").Value = 2.875
$ws.Range("M83").Value = 2.05
$ws.Range("N83").Value = 3.2
$ws.Range("O83").Value = 3.2
$ws.Range("P83").Value = -0.25
$ws.Range("Q83").Value = 1.825
$ws.Range("R83").Value = 1.975
$ws.Range("S83").Value = 2.5
$ws.Range("T83").Value = 1.825
$ws.Range("U83").Value = 1.975
$ws.Range("V83").Value = 1.05
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = 0.825
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.825
$ws.Range("AB83").Value = -1

# New values for row 84 (previously row 83's data)
$ws.Range("B84").Value = 7301364
$ws.Range("E84").Value = "Forge FC"
$ws.Range("F84").Value = "Atletico Ottawa"
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 1
$ws.Range("I84").Value = "A"
$ws.Range("J84").Value = 1.8
$ws.Range("K84").Value = 3.6
$ws.Range("L84").Value = 3.5
$ws.Range("M84").Value = 1.533
$ws.Range("N84").Value = 3.8
$ws.Range("O84").Value = 5
$ws.Range("P84").Value = -1
$ws.Range("Q84").Value = 1.975
$ws.Range("R84").Value = 1.825
$ws.Range("S84").Value = 2.5
$ws.Range("T84").Value = 1.9
$ws.Range("U84").Value = 1.9
$ws.Range("V84").Value = -1
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = 4
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 0.825
$ws.Range("AA84").Value = -1
$ws.Range("AB84").Value = 0.8999999999999999
